$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = -4.2598637128466095
$ws.Range("C2").Value = 4.397415285381463
$ws.Range("D2").Value = 1.2431418524331619
$ws.Range("E2").Value = 1.265501820836036

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 3.6266888763321674
$ws.Range("C3").Value = 5.7215295064976885
$ws.Range("D3").Value = 6.9050263606665823
$ws.Range("E3").Value = -0.94244729809990702

# Update selection to match final state
$ws.Range("B1:E3").Select()
